{"js": "// The commit turns the final \".\" (sometimes followed by a trailing\n// no-break space) of ten list-item sentences into a \";\" \u2014 the author\n// changed a run of standalone sentences into a semicolon-joined list.\n//\n// Paragraphs affected (unique leading text is enough to find each one):\n//   pStyle \"a0\" (6 items - AGM-AB contributions list)\n//   pStyle \"a\"  (4 items - VMI steps list: Access/Collection/Analysis/Logging)\nconst targets = [\n  \"Moreover, a novel AGM-AB algorithm is developed for detecting unknown malware functions from the benign program\",\n  \"Additionally, the African buffalo fitness module has been updated in the AGM manner to extract the features of the Virtual Machine Monitor (VMM)\",\n  \"Here, the introduced AGM-AB model investigates the guest operating system, system calls, and kernel data for classifying the malware and benign files\",\n  \"Also, the AGM-AB approach is tested by launching faults and malware functions to demonstrate the effectiveness of the AGM-AB method\",\n  \"Subsequently, the implementation of this proposed AGM-AB approach is done in the Python tool and the metrics are computed\",\n  \"Finally, the proposed method is evaluated by prevailing approaches in terms of recall, accuracy, AUC, FPR, precision, and F-measure\",\n  \"mechanism called COM/XPCOM that implements the VirtualBox API\",\n  \"Collection: It generates a memory dump of the Virtual Machine volatile memory\",\n  \"Analysis: It translates the low-level bytes into high-level information with the help of the Volatility tool, through the profile of the virtual machine and extracts objects from the operating system\",\n  \"Logging: It generates the log of the malware analysis\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const needle of targets) {\n  // Locate the paragraph that starts/contains this sentence.\n  let paragraph = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n      paragraph = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!paragraph) {\n    continue;\n  }\n\n  const fullRange = paragraph.getRange();\n  fullRange.load(\"text\");\n  await context.sync();\n\n  const text = fullRange.text;\n  // Trailing punctuation/whitespace to drop: an optional run of spaces /\n  // no-break spaces, the final \".\", then an optional trailing space / nbsp.\n  const match = /[\\u00A0 ]*\\.[\\u00A0 ]*$/.exec(text);\n  if (!match) {\n    continue;\n  }\n  const tail = text.substring(match.index);\n\n  // Isolate just that trailing chunk as its own sub-range so the rest of\n  // the paragraph (and its run formatting) is left untouched, then turn\n  // it into a single \";\".\n  const tailMatches = fullRange.search(tail, { matchCase: true });\n  tailMatches.load(\"text\");\n  await context.sync();\n\n  if (tailMatches.items.length === 0) {\n    continue;\n  }\n  const tailRange = tailMatches.items[tailMatches.items.length - 1];\n  tailRange.insertText(\";\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The commit turns the final \".\" (sometimes followed by a trailing\n# no-break space) of ten list-item sentences into a \";\" -- the author\n# changed a run of standalone sentences into a semicolon-joined list.\n#\n# Paragraphs affected (unique leading/inner text is enough to find each one):\n#   pStyle \"a0\" (6 items - AGM-AB contributions list)\n#   pStyle \"a\"  (4 items - VMI steps list: Access/Collection/Analysis/Logging)\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Moreover, a novel AGM-AB algorithm is developed for detecting unknown malware functions from the benign program\",\n  \"Additionally, the African buffalo fitness module has been updated in the AGM manner to extract the features of the Virtual Machine Monitor (VMM)\",\n  \"Here, the introduced AGM-AB model investigates the guest operating system, system calls, and kernel data for classifying the malware and benign files\",\n  \"Also, the AGM-AB approach is tested by launching faults and malware functions to demonstrate the effectiveness of the AGM-AB method\",\n  \"Subsequently, the implementation of this proposed AGM-AB approach is done in the Python tool and the metrics are computed\",\n  \"Finally, the proposed method is evaluated by prevailing approaches in terms of recall, accuracy, AUC, FPR, precision, and F-measure\",\n  \"mechanism called COM/XPCOM that implements the VirtualBox API\",\n  \"Collection: It generates a memory dump of the Virtual Machine volatile memory\",\n  \"Analysis: It translates the low-level bytes into high-level information with the help of the Volatility tool, through the profile of the virtual machine and extracts objects from the operating system\",\n  \"Logging: It generates the log of the malware analysis\"\n)\n\nforeach ($needle in $targets) {\n    foreach ($p in $d.Paragraphs) {\n        $full = $p.Range\n        $t = $full.Text\n        if ($t.Contains($needle)) {\n            # Paragraph.Range.Text ends with the paragraph mark (and,\n            # occasionally, a cell mark) -- strip those before looking for\n            # the trailing punctuation we want to swap for \";\".\n            $core = $t.TrimEnd([char]13, [char]7)\n            $m = [regex]::Match($core, \"[\\u00A0 ]*\\.[\\u00A0 ]*$\")\n            if ($m.Success) {\n                $tailStart = $full.Start + $m.Index\n                $tailEnd = $full.Start + $core.Length\n                $tailRange = $d.Range($tailStart, $tailEnd)\n                $tailRange.Text = \";\"\n            }\n            break\n        }\n    }\n}\n"}
